$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "IMAGE LINK" column (F) is being split into two columns:
#   F -> "CAT IMAGE LINK" (keeps the existing category-image values)
#   G -> new "HP IMAGE LINK" column with a placeholder banner image URL
$ws.Range("F1").Value = "CAT IMAGE LINK"

# Insert a brand new column at G (pushes GENERIC MESSAGE .. URL LINK one column right)
$ws.Columns("G").Insert()
$ws.Columns("G").ColumnWidth = 31.4
$ws.Range("G1").Value = "HP IMAGE LINK"
$ws.Range("G2").Value = "https://via.placeholder.com/500x50"
$ws.Range("G3").Value = "https://via.placeholder.com/500x50"

# Add a new trailing column "HP" (copy the formatting of the last existing
# column - URL LINK, now N - so the new column inherits matching styles)
$ws.Columns("N").Copy()
$ws.Columns("O").Insert()
$ws.Range("O1").Value = "HP"
$ws.Range("O2").Value = "YES"
$ws.Range("O3").Value = "NO"

# Update the view/selection state to match the edited workbook
$ws.Range("M18").Select()
$excel.ActiveWindow.ScrollColumn = 6
